$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("products")

$ws.Range("C2").Value = "f20adf46-e21b-4475-bc9c-1c157ef7610d"
$ws.Range("D2").Value = "fd858631-07a9-44f8-ab72-46928d93d15f"
$ws.Range("C3").Value = "f20adf46-e21b-4475-bc9c-1c157ef7610d"
$ws.Range("D3").Value = "fd858631-07a9-44f8-ab72-46928d93d15f"
$ws.Range("C4").Value = "f20adf46-e21b-4475-bc9c-1c157ef7610d"
$ws.Range("D4").Value = "fd858631-07a9-44f8-ab72-46928d93d15f"
$ws.Range("C5").Value = "f20adf46-e21b-4475-bc9c-1c157ef7610d"
$ws.Range("D5").Value = "fd858631-07a9-44f8-ab72-46928d93d15f"
$ws.Range("C6").Value = "f20adf46-e21b-4475-bc9c-1c157ef7610d"
$ws.Range("D6").Value = "dad2b61c-dcdc-460a-86a2-4b78aad14146"
$ws.Range("C7").Value = "f20adf46-e21b-4475-bc9c-1c157ef7610d"
$ws.Range("D7").Value = "dad2b61c-dcdc-460a-86a2-4b78aad14146"
$ws.Range("C8").Value = "f20adf46-e21b-4475-bc9c-1c157ef7610d"
$ws.Range("D8").Value = "dad2b61c-dcdc-460a-86a2-4b78aad14146"
$ws.Range("C9").Value = "f20adf46-e21b-4475-bc9c-1c157ef7610d"
$ws.Range("D9").Value = "dad2b61c-dcdc-460a-86a2-4b78aad14146"
$ws.Range("C10").Value = "f20adf46-e21b-4475-bc9c-1c157ef7610d"
$ws.Range("D10").Value = "dad2b61c-dcdc-460a-86a2-4b78aad14146"
$ws.Range("C11").Value = "f20adf46-e21b-4475-bc9c-1c157ef7610d"
$ws.Range("D11").Value = "e96b16a1-1db4-4fdb-9227-f63d14771655"
$ws.Range("C12").Value = "f20adf46-e21b-4475-bc9c-1c157ef7610d"
$ws.Range("D12").Value = "e96b16a1-1db4-4fdb-9227-f63d14771655"
$ws.Range("C13").Value = "f20adf46-e21b-4475-bc9c-1c157ef7610d"
$ws.Range("D13").Value = "e96b16a1-1db4-4fdb-9227-f63d14771655"
$ws.Range("C14").Value = "f20adf46-e21b-4475-bc9c-1c157ef7610d"
$ws.Range("D14").Value = "a8c771b0-910e-4b5b-bc0c-1902d4476642"
$ws.Range("C15").Value = "f20adf46-e21b-4475-bc9c-1c157ef7610d"
$ws.Range("D15").Value = "a8c771b0-910e-4b5b-bc0c-1902d4476642"
$ws.Range("C16").Value = "f20adf46-e21b-4475-bc9c-1c157ef7610d"
$ws.Range("D16").Value = "a8c771b0-910e-4b5b-bc0c-1902d4476642"
$ws.Range("C17").Value = "f20adf46-e21b-4475-bc9c-1c157ef7610d"
$ws.Range("D17").Value = "a8c771b0-910e-4b5b-bc0c-1902d4476642"
$ws.Range("C18").Value = "2e31f713-54ce-47d9-9a21-a14352ef6901"
$ws.Range("D18").Value = "2d78856f-75a7-49a8-9467-a8e65d9b2eae"
$ws.Range("C19").Value = "2e31f713-54ce-47d9-9a21-a14352ef6901"
$ws.Range("D19").Value = "2d78856f-75a7-49a8-9467-a8e65d9b2eae"
$ws.Range("C20").Value = "2e31f713-54ce-47d9-9a21-a14352ef6901"
$ws.Range("D20").Value = "2d78856f-75a7-49a8-9467-a8e65d9b2eae"
$ws.Range("C21").Value = "2e31f713-54ce-47d9-9a21-a14352ef6901"
$ws.Range("D21").Value = "2d78856f-75a7-49a8-9467-a8e65d9b2eae"
$ws.Range("C22").Value = "2e31f713-54ce-47d9-9a21-a14352ef6901"
$ws.Range("D22").Value = "0b138f75-7651-4e52-9618-73bd83bddad6"
$ws.Range("C23").Value = "2e31f713-54ce-47d9-9a21-a14352ef6901"
$ws.Range("D23").Value = "0b138f75-7651-4e52-9618-73bd83bddad6"
$ws.Range("C24").Value = "2e31f713-54ce-47d9-9a21-a14352ef6901"
$ws.Range("D24").Value = "0b138f75-7651-4e52-9618-73bd83bddad6"
$ws.Range("C25").Value = "2e31f713-54ce-47d9-9a21-a14352ef6901"
$ws.Range("D25").Value = "0a6649da-c760-4c59-a2f2-7ef3f5284103"
$ws.Range("C26").Value = "2e31f713-54ce-47d9-9a21-a14352ef6901"
$ws.Range("D26").Value = "0a6649da-c760-4c59-a2f2-7ef3f5284103"
$ws.Range("C27").Value = "2e31f713-54ce-47d9-9a21-a14352ef6901"
$ws.Range("D27").Value = "0a6649da-c760-4c59-a2f2-7ef3f5284103"
$ws.Range("C28").Value = "f81bc381-36e7-4364-940f-ec5aeac8cfdb"
$ws.Range("D28").Value = "17744f51-bf2e-4def-a6c5-dd29fd9f4a46"
$ws.Range("C29").Value = "f81bc381-36e7-4364-940f-ec5aeac8cfdb"
$ws.Range("D29").Value = "17744f51-bf2e-4def-a6c5-dd29fd9f4a46"
$ws.Range("C30").Value = "f81bc381-36e7-4364-940f-ec5aeac8cfdb"
$ws.Range("D30").Value = "17744f51-bf2e-4def-a6c5-dd29fd9f4a46"
$ws.Range("C31").Value = "f81bc381-36e7-4364-940f-ec5aeac8cfdb"
$ws.Range("D31").Value = "17744f51-bf2e-4def-a6c5-dd29fd9f4a46"
$ws.Range("C32").Value = "f81bc381-36e7-4364-940f-ec5aeac8cfdb"
$ws.Range("D32").Value = "17744f51-bf2e-4def-a6c5-dd29fd9f4a46"
$ws.Range("C33").Value = "f81bc381-36e7-4364-940f-ec5aeac8cfdb"
$ws.Range("D33").Value = "17744f51-bf2e-4def-a6c5-dd29fd9f4a46"
$ws.Range("C34").Value = "f81bc381-36e7-4364-940f-ec5aeac8cfdb"
$ws.Range("D34").Value = "17744f51-bf2e-4def-a6c5-dd29fd9f4a46"
$ws.Range("C35").Value = "f81bc381-36e7-4364-940f-ec5aeac8cfdb"
$ws.Range("D35").Value = "17744f51-bf2e-4def-a6c5-dd29fd9f4a46"
$ws.Range("C36").Value = "f81bc381-36e7-4364-940f-ec5aeac8cfdb"
$ws.Range("D36").Value = "17744f51-bf2e-4def-a6c5-dd29fd9f4a46"
$ws.Range("C37").Value = "f81bc381-36e7-4364-940f-ec5aeac8cfdb"
$ws.Range("D37").Value = "17744f51-bf2e-4def-a6c5-dd29fd9f4a46"
$ws.Range("C38").Value = "f81bc381-36e7-4364-940f-ec5aeac8cfdb"
$ws.Range("D38").Value = "a3d63bbd-beea-4c82-829f-df9289344c60"
$ws.Range("C39").Value = "f81bc381-36e7-4364-940f-ec5aeac8cfdb"
$ws.Range("D39").Value = "a3d63bbd-beea-4c82-829f-df9289344c60"
$ws.Range("C40").Value = "f81bc381-36e7-4364-940f-ec5aeac8cfdb"
$ws.Range("D40").Value = "a3d63bbd-beea-4c82-829f-df9289344c60"
$ws.Range("C41").Value = "f81bc381-36e7-4364-940f-ec5aeac8cfdb"
$ws.Range("D41").Value = "a3d63bbd-beea-4c82-829f-df9289344c60"
$ws.Range("C42").Value = "f81bc381-36e7-4364-940f-ec5aeac8cfdb"
$ws.Range("D42").Value = "a3d63bbd-beea-4c82-829f-df9289344c60"
$ws.Range("C43").Value = "f81bc381-36e7-4364-940f-ec5aeac8cfdb"
$ws.Range("D43").Value = "a3d63bbd-beea-4c82-829f-df9289344c60"
$ws.Range("C44").Value = "d92f6e1b-42b8-4244-8534-a90569cf1a89"
$ws.Range("D44").Value = "3807807b-b44a-4d4d-9bda-12646e61f6a3"
$ws.Range("C45").Value = "d92f6e1b-42b8-4244-8534-a90569cf1a89"
$ws.Range("D45").Value = "3807807b-b44a-4d4d-9bda-12646e61f6a3"
$ws.Range("C46").Value = "d92f6e1b-42b8-4244-8534-a90569cf1a89"
$ws.Range("D46").Value = "3807807b-b44a-4d4d-9bda-12646e61f6a3"
$ws.Range("C47").Value = "d92f6e1b-42b8-4244-8534-a90569cf1a89"
$ws.Range("D47").Value = "3807807b-b44a-4d4d-9bda-12646e61f6a3"
$ws.Range("C48").Value = "d92f6e1b-42b8-4244-8534-a90569cf1a89"
$ws.Range("D48").Value = "3807807b-b44a-4d4d-9bda-12646e61f6a3"
$ws.Range("C49").Value = "d92f6e1b-42b8-4244-8534-a90569cf1a89"
$ws.Range("D49").Value = "3807807b-b44a-4d4d-9bda-12646e61f6a3"
$ws.Range("C50").Value = "d92f6e1b-42b8-4244-8534-a90569cf1a89"
$ws.Range("D50").Value = "3807807b-b44a-4d4d-9bda-12646e61f6a3"
$ws.Range("C51").Value = "d92f6e1b-42b8-4244-8534-a90569cf1a89"
$ws.Range("D51").Value = "3807807b-b44a-4d4d-9bda-12646e61f6a3"
$ws.Range("C52").Value = "d92f6e1b-42b8-4244-8534-a90569cf1a89"
$ws.Range("D52").Value = "3807807b-b44a-4d4d-9bda-12646e61f6a3"
$ws.Range("C53").Value = "d92f6e1b-42b8-4244-8534-a90569cf1a89"
$ws.Range("D53").Value = "437d4659-7f77-4f1c-be9b-5589924e665a"
$ws.Range("C54").Value = "d92f6e1b-42b8-4244-8534-a90569cf1a89"
$ws.Range("D54").Value = "437d4659-7f77-4f1c-be9b-5589924e665a"
$ws.Range("C55").Value = "d92f6e1b-42b8-4244-8534-a90569cf1a89"
$ws.Range("D55").Value = "437d4659-7f77-4f1c-be9b-5589924e665a"
$ws.Range("C56").Value = "d92f6e1b-42b8-4244-8534-a90569cf1a89"
$ws.Range("D56").Value = "437d4659-7f77-4f1c-be9b-5589924e665a"
$ws.Range("C57").Value = "d92f6e1b-42b8-4244-8534-a90569cf1a89"
$ws.Range("D57").Value = "437d4659-7f77-4f1c-be9b-5589924e665a"
